# feat: add 2022-Q1 data
#
# The workbook gains a new "2022-Q1" fund-holding detail sheet (placed
# between the existing "2021-Q4" sheet and "总计"), and the "总计" summary
# sheet gains a new top data row for "2022-Q1" (the old "2021-Q4" summary
# row shifts down).
#
# The original "总计" worksheet is the template for the new "2022-Q1"
# sheet: it is duplicated first (the duplicate becomes the refreshed
# "总计"), then the original sheet object itself is renamed to "2022-Q1"
# and its header/data are rewritten - which keeps its original formatting
# (borders/bold/centering) for the new sheet, matching how the source
# data is laid out.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$totalSrc = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Clone "总计" -> the clone will become the refreshed "总计" sheet;
#    it lands right after the original (i.e. after "2021-Q4").
# ---------------------------------------------------------------------
$totalSrc.Copy($null, $totalSrc) | Out-Null
$totalNew = $wb.Worksheets.Item($totalSrc.Index + 1)
$totalNew.Name = "总计 (tmp)"

# ---------------------------------------------------------------------
# 2. Turn the original "总计" sheet into the "2022-Q1" detail sheet.
# ---------------------------------------------------------------------
$totalSrc.Name = "2022-Q1"
$q1 = $totalSrc

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# New header cells E1:H1 need the same header styling as B1:D1 - copy it
# across from the (already correctly styled) D1.
$q1.Range("D1").Copy() | Out-Null
$q1.Range("E1:H1").PasteSpecial(-4122) | Out-Null

# Fund-holding data row. B2 ("005706") and D2:G2 are text in the source
# data (e.g. "2.37") rather than numbers - force text formatting before
# assigning so the literal strings (leading zero, fixed decimals) stick,
# then drop the number-format override again so no stray style lingers.
$textCells = $q1.Range("B2:B2,D2:G2")
$textCells.NumberFormat = "@"
$q1.Range("B2").Value = "005706"
$q1.Range("C2").Value = "兴业龙腾双益平衡混合"
$q1.Range("D2").Value = "2.37"
$q1.Range("E2").Value = "27.85"
$q1.Range("F2").Value = "4.33"
$q1.Range("G2").Value = "0.1026"
$textCells.ClearFormats() | Out-Null

# H2 (仓位排名) stays a real number.
$q1.Range("H2").Value = 2

# ---------------------------------------------------------------------
# 3. Prepend the "2022-Q1" row to the refreshed "总计" sheet.
# ---------------------------------------------------------------------
$totalNew.Name = "总计"
$total = $totalNew

$total.Rows.Item(2).Insert() | Out-Null

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.1
$total.Range("B2:D2").ClearFormats() | Out-Null

# Match the row-index styling used by every other data row (copy the
# header's style onto the new index cell).
$total.Range("B1").Copy() | Out-Null
$total.Range("A2").PasteSpecial(-4122) | Out-Null

# The old "2021-Q4" row shifted from row 2 to row 3; fix its index column.
$total.Range("A3").Value = 1

# ---------------------------------------------------------------------
# Restore the original active sheet/selection.
# ---------------------------------------------------------------------
$q4.Activate() | Out-Null
$q4.Range("A1").Select() | Out-Null
